# The source diff reshuffles the data records that sit in rows 2-7 of the
# "Artfynd" sheet into new row positions (same six records, same columns,
# just relocated - e.g. the record that was on row 2 ends up on row 4, the
# one on row 7 ends up on row 2, and so on). Column I ("Antal") stores
# numeral-looking values as text in this workbook, so it is written with a
# leading apostrophe to keep Excel from coercing it back to a number.
#
# Only cells whose value actually differs between the old and new record at
# a given row are written, to avoid touching cells that happen to already
# hold the correct value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column I to text storage so the numeral strings round-trip as text.
$ws.Range("I2:I7").NumberFormat = "@"

$ws.Range("A2").Value = 111543968
$ws.Range("B2").Value = 57487
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 205998
$ws.Range("F2").Value = "Nordfladdermus"
$ws.Range("G2").Value = "Eptesicus nilssonii"
$ws.Range("H2").Value = "(A.Keyserling & Blasius, 1839)"
$ws.Range("I2").Value = "256"
$ws.Range("J2").Value = ""

$ws.Range("A3").Value = 111545328
$ws.Range("B3").Value = 57494
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 205992
$ws.Range("F3").Value = "Vattenfladdermus"
$ws.Range("G3").Value = "Myotis daubentonii"
$ws.Range("H3").Value = "(Kuhl, 1817)"
$ws.Range("I3").Value = "1"

$ws.Range("A4").Value = 111543957
$ws.Range("I4").Value = "1"
$ws.Range("P4").Value = "Orsa Viborg, glänta i skogsparti, Dlr"
$ws.Range("Q4").Value = 480406.6045043401
$ws.Range("R4").Value = 6772745.04339793

$ws.Range("A5").Value = 111545414
$ws.Range("I5").Value = "9"
$ws.Range("J5").Value = "registreringar"
$ws.Range("P5").Value = "Orsa Viborg, glänta i mitten av skogsparti, Dlr"
$ws.Range("Q5").Value = 480487.2503558649
$ws.Range("R5").Value = 6772784.264016891

$ws.Range("A6").Value = 111545323
$ws.Range("I6").Value = "2"
$ws.Range("J6").Value = ""
$ws.Range("P6").Value = "Orsa Viborg, intill en grupp med hålträd, Dlr"
$ws.Range("Q6").Value = 480427.8053356989
$ws.Range("R6").Value = 6772811.198980245

$ws.Range("A7").Value = 111545401
$ws.Range("I7").Value = "6"
$ws.Range("J7").Value = "registreringar"
$ws.Range("P7").Value = "Orsa Viborg, glänta i mitten av skogsparti, Dlr"
$ws.Range("Q7").Value = 480487.2503558649
$ws.Range("R7").Value = 6772784.264016891
